$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center the "Complete" header in D3 (keeps its bold font, adds center alignment)
$ws.Range("D3").HorizontalAlignment = -4108  # xlCenter

# Mark milestone items as complete with an "X", centered
$ws.Range("D4").Value = "X"
$ws.Range("D4").HorizontalAlignment = -4108  # xlCenter

$ws.Range("D12").Value = "X"
$ws.Range("D12").HorizontalAlignment = -4108  # xlCenter

# Update the view/selection state to match what was captured on save
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D13").Select()
